# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# Row -> new value updates for "展览" sheet
$exhibitUpdates = @{
    2  = 590
    4  = 1270
    5  = 1117
    6  = 14144
    7  = 15774
    9  = 61
    11 = 193
    17 = 31
    18 = 82
    20 = 1223
    22 = 67
    23 = 15
    24 = 6227
    26 = 1096
    27 = 5583
    28 = 76
    30 = 125
    31 = 4571
    32 = 5
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new value updates for "全部类型" sheet
$allUpdates = @{
    2  = 590
    4  = 1270
    5  = 1117
    6  = 14144
    7  = 15774
    9  = 61
    11 = 193
    17 = 31
    18 = 82
    20 = 1223
    22 = 67
    24 = 15
    25 = 6227
    27 = 1096
    28 = 5583
    29 = 76
    31 = 125
    32 = 4571
    33 = 5
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
